$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new scheduled posts to the bottom of the table
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Ground Squirrel"
$ws.Range("D21").Value = "Shot"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Lady Quail"
$ws.Range("D22").Value = "Shot"

# Move the active selection as recorded in the saved workbook
$ws.Range("N3").Select()
